# AB#10856 - Add Immunization Recommendations section to the report template
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "Immunizations"
$ws.Name = "Immunizations"

# Insert a new title row above the existing "Immunization" history table
$ws.Rows.Item(1).Insert()

# New title cell for the (now shifted-down) history table
$ws.Range("A1").Value = "Immunization History"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 16
$ws.Rows.Item(1).RowHeight = 21

# New "Immunization Recommendations" section, two blank rows below the
# existing table (which now occupies rows 2-4)
$ws.Range("A7").Value = "Immunization Recommendations"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Size = 16
$ws.Rows.Item(7).RowHeight = 21

# Header row for the recommendations table
$ws.Range("A8").Value = "Immunization"
$ws.Range("B8").Value = "Due Date"
$ws.Range("C8").Value = "Status"
$ws.Range("A8:C8").Font.Bold = $true

# Template placeholder rows for the recommendations table body
$ws.Range("A9").Value = "{d.recommendations[i].immunization}"
$ws.Range("B9").Value = "{d.recommendations[i].due_date}"
$ws.Range("C9").Value = "{d.recommendations[i].status}"

$ws.Range("A10").Value = "{d.recommendations[i+1].immunization}"
$ws.Range("B10").Value = "{d.recommendations[i+1].due_date}"
$ws.Range("C10").Value = "{d.recommendations[i+1].status}"
